# Actualización automática 2025-09-08 09:31:10
# Set the PRESUPUESTO (column G) values to 0 for several clients on the
# "VENTA MENSUAL" sheet, and update the totals row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 0
